$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "404÷6=" "185÷5="
Replace-Text "548÷3=" "449÷9="
Replace-Text "303÷6=" "840÷9="
Replace-Text "921÷9=" "505÷8="
Replace-Text "712÷7=" "289÷2="
Replace-Text "908÷8=" "762÷6="
Replace-Text "777÷3=" "643÷8="
Replace-Text "610÷7=" "322÷7="
Replace-Text "408÷2=" "597÷2="
Replace-Text "532÷3=" "251÷7="
Replace-Text "880÷4=" "427÷3="
Replace-Text "652÷6=" "708÷3="
Replace-Text "695÷2=" "373÷6="
Replace-Text "614÷7=" "442÷3="
Replace-Text "469÷2=" "330÷2="
Replace-Text "837÷8=" "844÷4="
Replace-Text "422÷3=" "202÷9="
Replace-Text "448÷2=" "641÷5="
Replace-Text "805÷6=" "824÷5="
Replace-Text "788÷9=" "226÷6="
Replace-Text "971÷7=" "173÷2="
Replace-Text "672÷7=" "774÷6="
Replace-Text "889÷5=" "785÷4="
Replace-Text "395÷8=" "845÷3="
Replace-Text "649÷8=" "683÷4="
